$wb = $excel.ActiveWorkbook

# Sheet "INPFC" header B1: "Latitude (upper limit)" -> "northlimit_latitude"
$ws1 = $wb.Worksheets.Item("INPFC")
$ws1.Range("B1").Value = "northlimit_latitude"

# Update selection on INPFC sheet from B3 to B2
$ws1.Activate()
$ws1.Range("B2").Select()
